$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.184.08'
$ws.Range("E2").Value = '  -1.47%  '
$ws.Range("D3").Value = '3.493.73'
$ws.Range("E3").Value = '  -3.73%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '198.98'
$ws.Range("E5").Value = '  +1.84%  '
$ws.Range("D6").Value = '550.12'
$ws.Range("E6").Value = '  -4.89%  '
$ws.Range("D7").Value = '3.490.77'
$ws.Range("E7").Value = '  -3.69%  '
$ws.Range("D8").Value = '0.604'
$ws.Range("E8").Value = '  -2.66%  '
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  -0.31%  '
$ws.Range("D10").Value = '0.650'
$ws.Range("E10").Value = '  -4.10%  '
$ws.Range("D11").Value = '62.47'
$ws.Range("E11").Value = '  +12.09%  '
$ws.Range("D12").Value = '0.143'
$ws.Range("E12").Value = '  -6.98%  '
$ws.Range("D13").Value = '0.0000268'
$ws.Range("E13").Value = '  -9.31%  '
$ws.Range("D14").Value = '9.76'
$ws.Range("E14").Value = '  -3.52%  '
$ws.Range("D15").Value = '4.043.23'
$ws.Range("E15").Value = '  -3.92%  '
$ws.Range("D16").Value = '3.487.82'
$ws.Range("E16").Value = '  -3.99%  '
$ws.Range("D17").Value = '0.123'
$ws.Range("E17").Value = '  -2.26%  '
$ws.Range("D18").Value = '66.830.75'
$ws.Range("E18").Value = '  -1.93%  '
$ws.Range("D19").Value = '18.23'
$ws.Range("E19").Value = '  -1.60%  '
$ws.Range("D20").Value = '11.75'
$ws.Range("E20").Value = '  -6.34%  '
$ws.Range("D21").Value = '1.02'
$ws.Range("E21").Value = '  -5.75%  '
$ws.Range("D22").Value = '387.81'
$ws.Range("E22").Value = '  -3.85%  '
$ws.Range("D23").Value = '3.98'
$ws.Range("E23").Value = '  -5.68%  '
$ws.Range("D24").Value = '11.85'
$ws.Range("E24").Value = '  -7.41%  '
$ws.Range("D25").Value = '82.13'
$ws.Range("E25").Value = '  -4.55%  '
$ws.Range("D26").Value = '3.83'
$ws.Range("E26").Value = '  -0.89%  '
$ws.Range("D27").Value = '12.12'
$ws.Range("E27").Value = '  -3.97%  '
$ws.Range("D28").Value = '2.78'
$ws.Range("E28").Value = '  -5.70%  '
$ws.Range("D29").Value = '8.74'
$ws.Range("E29").Value = '  -4.40%  '
$ws.Range("D30").Value = '30.91'
$ws.Range("E30").Value = '  -2.62%  '
$ws.Range("D31").Value = '674.39'
$ws.Range("E31").Value = '  -2.02%  '
$ws.Range("D32").Value = '6.93'
$ws.Range("E32").Value = '  -14.69%  '
$ws.Range("D33").Value = '11.65'
$ws.Range("E33").Value = '  -4.70%  '
$ws.Range("D34").Value = '63.40'
$ws.Range("E34").Value = '  -2.18%  '
$ws.Range("D35").Value = '0.109'
$ws.Range("E35").Value = '  -7.44%  '
$ws.Range("D36").Value = '38.16'
$ws.Range("E36").Value = '  -10.44%  '
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.32%  '
$ws.Range("D38").Value = '0.396'
$ws.Range("E38").Value = '  -5.12%  '
$ws.Range("D39").Value = '0.997'
$ws.Range("E39").Value = '  -0.23%  '
$ws.Range("D40").Value = '3.058.52'
$ws.Range("E40").Value = '  -3.03%  '
$ws.Range("D41").Value = '0.129'
$ws.Range("E41").Value = '  -5.08%  '
$ws.Range("D42").Value = '2.97'
$ws.Range("E42").Value = '  -5.46%  '
$ws.Range("D43").Value = '0.0₃0670'
$ws.Range("E43").Value = '  -16.08%  '
$ws.Range("D44").Value = '2.76'
$ws.Range("E44").Value = '  +6.56%  '
$ws.Range("D45").Value = '2.49'
$ws.Range("E45").Value = '  -13.44%  '
$ws.Range("D46").Value = '2.73'
$ws.Range("E46").Value = '  -6.12%  '
$ws.Range("D47").Value = '0.0393'
$ws.Range("E47").Value = '  -7.28%  '
$ws.Range("D48").Value = '0.126'
$ws.Range("E48").Value = '  -5.05%  '
$ws.Range("D49").Value = '136.77'
$ws.Range("E49").Value = '  -4.05%  '
$ws.Range("D50").Value = '2.89'
$ws.Range("E50").Value = '  -6.39%  '
$ws.Range("D51").Value = '8.16'
$ws.Range("E51").Value = '  -7.83%  '
